$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student Points")

# ---------------------------------------------------------------------------
# Group A, Community Site
# Insert a new rubric line ("A page showing all messages") right before the
# existing "Database works locally" row (original row 13), and drop that
# row's score from 15 -> 10 (the 5 points moved to the new line item).
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value = "A page showing all messages"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 5
$ws.Range("B14").Value = 10
$ws.Range("C14").Value = 10

# ---------------------------------------------------------------------------
# Group B, Fan Site
# Same pattern: "A page showing all stories" before "Database works locally"
# (original row 35, now row 36 after the Group A insert above).
# ---------------------------------------------------------------------------
$ws.Rows.Item(36).Insert()
$ws.Range("A36").Value = "A page showing all stories"
$ws.Range("B36").Value = 5
$ws.Range("C36").Value = 5
$ws.Range("B37").Value = 10
$ws.Range("C37").Value = 10

# ---------------------------------------------------------------------------
# Group C, Info Site
# The comment on this group changes, the Exercise grade drops to 0, and the
# usual new rubric line ("A page showing all comments") is inserted before
# "Database works locally" -- but this group only gets partial credit, with
# grader comments in column D explaining the shortfalls.
# ---------------------------------------------------------------------------
$ws.Rows.Item(59).Insert()
$ws.Range("A59").Value = "A page showing all comments"
$ws.Range("B59").Value = 5
$ws.Range("C59").Value = 3
$ws.Range("D59").Value = "Not found"

$ws.Range("B60").Value = 10
$ws.Range("C60").Value = 10
$ws.Range("D60").Value = "Comment is stored"

$ws.Range("B61").Value = 10
$ws.Range("C61").Value = 5
$ws.Range("D61").Value = "Can't check"

$ws.Range("C53").Value = 0
$ws.Range("D53").Value = "Not done"

$ws.Range("A50").Value = "Good start."

# ---------------------------------------------------------------------------
# View state to roughly match the saved workbook (active sheet scrolled down
# to the newly-edited Group C block).
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("A49:D69").Select()
